$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (32) with the final entry for the project deliverable.
$ws.Range("B32").Value = "Gaizka"
$ws.Range("C32").Value = "Esquema entidad relación"

# D32 gets the same date value/style as D31 (date stored as serial 45784, format like "7-may").
$ws.Range("D32").Value = $ws.Range("D31").Value2
$ws.Range("D32").NumberFormat = $ws.Range("D31").NumberFormat

# Move the active selection like in the target workbook (selection moved to D33).
$ws.Activate()
$ws.Range("D33").Select()
